$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = "CO"

# AC2 (PrimaryInsurance_ZipCode) is stored as text (it was "303740800" in
# the shared-string table, not a number). Temporarily mark the cell as
# Text so Excel keeps the new all-digit value as a string instead of
# auto-converting it to a number, then restore the cell's original style.
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = "802175747"
$ws.Range("AC2").Style = "Normal"
